# Re-simulated Week 17, factoring in more player injuries
# M.Sanders is ruled out for the week; his row is removed from both sheets
# and the rest of the backfield / receiving corps gets freshly simulated
# stat lines (with rows below shifting up to fill the gap on Rushing).

$wb = $excel.ActiveWorkbook
$rushing   = $wb.Worksheets.Item("Rushing")
$receiving = $wb.Worksheets.Item("Receiving")

# ---------------------------------------------------------------------------
# Rushing sheet
# ---------------------------------------------------------------------------
# Row 4 was M.Sanders - now B.Scott's line (player index column A unchanged)
$rushing.Range("B4").Value = "B.Scott"
$rushing.Range("C4").Value = 118
$rushing.Range("D4").Value = 83
$rushing.Range("E4").Value = 13
$rushing.Range("F4").Value = 30

# Row 5 was B.Scott - now J.Howard's line
$rushing.Range("B5").Value = "J.Howard"
$rushing.Range("C5").Value = 50
$rushing.Range("D5").Value = 17
$rushing.Range("E5").Value = 13
$rushing.Range("F5").Value = 15

# Row 6 stays K.Gainwell, with a new (much quieter) stat line
$rushing.Range("C6").Value = 6
$rushing.Range("D6").Value = 5
$rushing.Range("F6").Value = 5

# Row 7 was J.Howard - becomes J.Reagor's line, shifted up from row 8, and
# the player-index column A advances since M.Sanders' row is gone
$rushing.Range("A7").Value = 6
$rushing.Range("B7").Value = "J.Reagor"
$rushing.Range("C7").Value = 5
$rushing.Range("D7").Value = 1
$rushing.Range("E7").Value = 0
$rushing.Range("F7").Value = 0

# Row 8 was J.Reagor - becomes Q.Watkins' line, shifted up from row 9
$rushing.Range("A8").Value = 7
$rushing.Range("B8").Value = "Q.Watkins"
$rushing.Range("C8").Value = 1
$rushing.Range("D8").Value = 0
$rushing.Range("E8").Value = 0
$rushing.Range("F8").Value = 0

# Row 9 was Q.Watkins - becomes M.Walker's line, shifted up from row 10
$rushing.Range("A9").Value = 8
$rushing.Range("B9").Value = "M.Walker"
$rushing.Range("C9").Value = 0
$rushing.Range("D9").Value = 0
$rushing.Range("E9").Value = 1
$rushing.Range("F9").Value = 0

# Old row 10 (M.Walker's original line) is now redundant - drop it so the
# sheet ends at row 9.
$rushing.Rows.Item(10).Delete()

# ---------------------------------------------------------------------------
# Receiving sheet - same injury shuffle, but no rows are removed here; only
# the top three receiving backs get new names / numbers.
# ---------------------------------------------------------------------------
$receiving.Range("B2").Value = "J.Howard"

$receiving.Range("B3").Value = "K.Gainwell"
$receiving.Range("C3").Value = 23

$receiving.Range("B4").Value = "B.Scott"
$receiving.Range("C4").Value = 30

# ---------------------------------------------------------------------------
# View state: Receiving becomes the active tab with E3 selected; Rushing
# keeps a selection over its now-shorter data range.
# ---------------------------------------------------------------------------
$rushing.Activate()
$rushing.Range("A7:F9").Select()

$receiving.Activate()
$receiving.Range("E3").Select()
